# Corrected excel sheets for application fix issues
$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("A4").Value = 100
$wsSummary.Range("B4").Value = 100
$wsSummary.Activate()
$wsSummary.Range("B4").Select()

# --- Repayment Schedule sheet ---
$wsRepay = $wb.Worksheets.Item("Repayment Schedule")
$wsRepay.Range("I2").Value = 100
$wsRepay.Range("K2").Value = 100
$wsRepay.Range("L2").Value = 100
$wsRepay.Activate()
$wsRepay.Range("L2").Select()

# --- Transactions sheet ---
$wsTrans = $wb.Worksheets.Item("Transactions")
$wsTrans.Range("E2").Value = 100
$wsTrans.Range("H2").Value = 100
$wsTrans.Activate()
$wsTrans.Range("I11").Select()
